$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_data3")

# Clear the old table range (rows 14-22, cols B:C)
$ws.Range("B14:C22").Clear()

# Table 1: Sweden_Pay_Now_Direct_debit header + 4 key/value rows (rows 16-20)
$ws.Range("B16").Value = "Sweden_Pay_Now_Direct_debit"
$ws.Range("B17").Value = "key4"
$ws.Range("C17").Value = "value4"
$ws.Range("B18").Value = "key4"
$ws.Range("C18").Value = "value4"
$ws.Range("B19").Value = "key4"
$ws.Range("C19").Value = "value4"
$ws.Range("B20").Value = "key4"
$ws.Range("C20").Value = "value4"

# Row 21 intentionally left blank as a spacer between the two tables

# Table 2: Sweden_Pay_Now_Card header + 4 key/value rows (rows 22-26)
$ws.Range("B22").Value = "Sweden_Pay_Now_Card"
$ws.Range("B23").Value = "key4"
$ws.Range("C23").Value = "value4"
$ws.Range("B24").Value = "key4"
$ws.Range("C24").Value = "value4"
$ws.Range("B25").Value = "key4"
$ws.Range("C25").Value = "value4"
$ws.Range("B26").Value = "key4"
$ws.Range("C26").Value = "value4"
